$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (e.g. H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill I2:I36 with a constant 1, and J2:J36 with the same value as the
# corresponding H column cell on that row.
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $hval = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 10).Value = $hval
}
